$d = $word.ActiveDocument
$tbl = $d.Tables(1)
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# Change 1: "Eve Portal" -> "Arknights" (wrapped in proofErr spellStart/End)
# ---------------------------------------------------------------------------
$cell = $tbl.Cell(1, 3)
$p = $cell.Range.Paragraphs(1)
$xml = "<w:p $wns><w:pPr><w:rPr><w:lang w:val=`"en-US`"/></w:rPr></w:pPr>" +
       "<w:proofErr w:type=`"spellStart`"/>" +
       "<w:r><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>Arknights</w:t></w:r>" +
       "<w:proofErr w:type=`"spellEnd`"/></w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# Change 2: "Конфиденциальная информация" cell content rewritten
# ---------------------------------------------------------------------------
$cell = $tbl.Cell(5, 3)
$p = $cell.Range.Paragraphs(1)
$xml = "<w:p $wns>" +
       "<w:r><w:t xml:space=`"preserve`">Информация о соц. сетях пользователя, </w:t></w:r>" +
       "<w:proofErr w:type=`"spellStart`"/>" +
       "<w:r><w:t>токены</w:t></w:r>" +
       "<w:proofErr w:type=`"spellEnd`"/>" +
       "<w:r><w:t xml:space=`"preserve`"> доступа и сессии, информация об аккаунте и средствах на аккаунте пользователя</w:t></w:r>" +
       "</w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# Change 3: move the "_GoBack" bookmark from the 1.1.4 result cell down to
# the (newly populated) 1.7 result cell, and fill in the empty result cells
# for requirements 1.2 - 1.7.
# ---------------------------------------------------------------------------

# 1.1.4 result cell (row 12) keeps its text but loses the bookmark.
$cell = $tbl.Cell(12, 3)
$p = $cell.Range.Paragraphs(1)
$xml = "<w:p $wns>" +
       "<w:r><w:rPr><w:lang w:val=`"be-BY`"/></w:rPr><w:t xml:space=`"preserve`">Флаг </w:t></w:r>" +
       "<w:proofErr w:type=`"spellStart`"/>" +
       "<w:r><w:rPr><w:lang w:val=`"de-DE`"/></w:rPr><w:t>debuggable</w:t></w:r>" +
       "<w:proofErr w:type=`"spellEnd`"/>" +
       "<w:r><w:rPr><w:lang w:val=`"de-DE`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r>" +
       "<w:r><w:t xml:space=`"preserve`">установлен в </w:t></w:r>" +
       "<w:proofErr w:type=`"spellStart`"/>" +
       "<w:r><w:rPr><w:lang w:val=`"de-DE`"/></w:rPr><w:t>false</w:t></w:r>" +
       "<w:proofErr w:type=`"spellEnd`"/>" +
       "<w:r><w:rPr><w:lang w:val=`"de-DE`"/></w:rPr><w:t>.</w:t></w:r>" +
       "<w:r><w:t xml:space=`"preserve`"> Секретных ключей не обнаружено.</w:t></w:r>" +
       "</w:p>"
$p.Range.InsertXML($xml)

# 1.2 result cell (row 14): "Отсутствуют"
$cell = $tbl.Cell(14, 3)
$p = $cell.Range.Paragraphs(1)
$xml = "<w:p $wns><w:r><w:t>Отсутствуют</w:t></w:r></w:p>"
$p.Range.InsertXML($xml)

# 1.3 result cell (row 15): "Отсутсвует" (spelling kept as-authored, flagged)
$cell = $tbl.Cell(15, 3)
$p = $cell.Range.Paragraphs(1)
$xml = "<w:p $wns><w:proofErr w:type=`"spellStart`"/><w:r><w:t>Отсутсвует</w:t></w:r><w:proofErr w:type=`"spellEnd`"/></w:p>"
$p.Range.InsertXML($xml)

# 1.4 result cell (row 16): "Отсутствуют"
$cell = $tbl.Cell(16, 3)
$p = $cell.Range.Paragraphs(1)
$xml = "<w:p $wns><w:r><w:t>Отсутствуют</w:t></w:r></w:p>"
$p.Range.InsertXML($xml)

# 1.5 result cell (row 17): "Не замечено каких-либо нарушений и уязвимостей"
$cell = $tbl.Cell(17, 3)
$p = $cell.Range.Paragraphs(1)
$xml = "<w:p $wns><w:r><w:t>Не замечено каких-либо нарушений и уязвимостей</w:t></w:r></w:p>"
$p.Range.InsertXML($xml)

# 1.6 result cell (row 18): "Не требуется"
$cell = $tbl.Cell(18, 3)
$p = $cell.Range.Paragraphs(1)
$xml = "<w:p $wns><w:r><w:t>Не требуется</w:t></w:r></w:p>"
$p.Range.InsertXML($xml)

# 1.7 result cell (row 19): new text, and the "_GoBack" bookmark now lives here
$cell = $tbl.Cell(19, 3)
$p = $cell.Range.Paragraphs(1)
$xml = "<w:p $wns>" +
       "<w:r><w:t>Все конфиденциальные данные успешно утекли в созданную резервную копию.</w:t></w:r>" +
       "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/>" +
       "</w:p>"
$p.Range.InsertXML($xml)

Write-Host "done"
